$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 2 (shifts existing rows 2..162 down to 3..163)
$ws.Rows("2:2").Insert()

# Fill in the new row 2 with the "tst" test entry
$ws.Range("A2").Value = 1111
$ws.Range("B2").Value = "tst"
$ws.Range("C2").Value = "3553-Vacaciones"

# Copy the date-formatted style from row 3 (the old row 2, shifted down) onto D2:E2
# before writing the date values, so the new cells reuse the existing date style
# instead of Excel minting a brand-new number-format style.
$ws.Range("D3:E3").Copy()
$ws.Range("D2:E2").PasteSpecial(-4122)
$ws.Range("D2").Value = "12/26/2023"
$ws.Range("E2").Value = "1/8/2024"

$ws.Range("F2").Value = 6
$ws.Range("G2").Value = 14

# The insert above shifted the whole sheet (including the block of blank
# placeholder rows at the bottom) down by one row. Only the data block
# (rows 2-30) should have grown though - the long tail of blank rows must
# stay exactly where it was, so remove one now-duplicate blank row to
# compensate.
$ws.Rows("31:31").Delete()

# Re-extend the blank placeholder tail by one row so the sheet still ends
# at row 163, copying the formatting (short-date style on D/E) from the
# last existing placeholder row.
$ws.Range("D162:E162").Copy()
$ws.Range("D163:E163").PasteSpecial(-4122)

# Restore the cursor/selection to where the editor left it.
[void]$ws.Range("B4").Select()
